# Autocomplete para veículo e inclusão de novos ícones para mapa.
# This script reproduces, via Excel COM-interop, the edits captured in the
# diff: new "Plan3" sheet content (rota/hora/tipo/calendário table), the
# sheet becoming the active tab/selection, and sheet1 losing its
# "tabSelected" flag as a consequence of Plan3 becoming active.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Plan3")

# --- Header row -------------------------------------------------------
# Writing cell-by-cell in this precise order reproduces the shared-string
# table ordering seen in the target workbook (rota, rota 001, hora ini,
# hora fim, tipo, data, calendário, dia, ano letivo).
$ws3.Range("A1").Value = "rota"
$ws3.Range("A2").Value = "rota 001"

$ws3.Range("B1").Value = "hora ini"
$ws3.Range("C1").Value = "hora fim"

$ws3.Range("D1").Value = "tipo"

$ws3.Range("F1").Value = "data"

$ws3.Range("D2").Value = "calendário"
$ws3.Range("D3").Value = "dia"

$ws3.Range("E1").Value = "ano letivo"

# --- Data rows ----------------------------------------------------------
$ws3.Range("B2").Value = 0.29166666666666669
$ws3.Range("B2").NumberFormat = "h:mm"

$ws3.Range("C2").Value = 0.33333333333333331
$ws3.Range("C2").NumberFormat = "h:mm"

$ws3.Range("E2").Value = 2013

$ws3.Range("B3").Value = 0.45833333333333331
$ws3.Range("B3").NumberFormat = "h:mm"

$ws3.Range("C3").Value = 0.5
$ws3.Range("C3").NumberFormat = "h:mm"

$ws3.Range("F3").Value = 41397
$ws3.Range("F3").NumberFormat = "m/d/yy"

# --- Column widths (best-fit like the original authoring session) -------
$ws3.Columns.Item(4).ColumnWidth = 9.43
$ws3.Columns.Item(5).ColumnWidth = 9.43
$ws3.Columns.Item(6).ColumnWidth = 9.75

# --- Make Plan3 the active sheet/tab and set its selection ---------------
$ws3.Activate()
$ws3.Range("G3").Select()
